$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Heap Report from Test")

# Add a new log entry (row 11) documenting the "Reserve_Stack_Space" branch test
$ws.Range("A11").NumberFormat = "[$-409]m/d/yy\ h:mm\ AM/PM;@"
$ws.Range("A11").Value = 43426.651388888888
$ws.Range("B11").Value = "PC"
$ws.Range("C11").Value = "Release"
$ws.Range("D11").Value = "Reserve_Stack_Space"

$ws.Range("F11").NumberFormat = "#,##0"
$ws.Range("F11").WrapText = $true
$ws.Range("F11").Value = "7,916`n7,916`n7,916"

$ws.Range("G11").NumberFormat = "#,##0"
$ws.Range("G11").WrapText = $true
$ws.Range("G11").Value = "42,292`n42,132`n42,372"

$ws.Range("H11").NumberFormat = "#,##0"
$ws.Range("I11").WrapText = $true

$ws.Rows.Item(11).RowHeight = 43.2

# Expand the log table to include the new row
$ws.ListObjects.Item("Table1").Resize($ws.Range("A1:J11"))

$wb.Save()
